# Applies the cryptos.xlsx data refresh described in the commit diff:
# updates Price (D) / Volume(1h) (E) figures for existing rows, and
# inserts a new "EnergySwap" row at position 49, pushing FraxShare and
# ordi down one row and dropping the old MultiversX row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text reads like a plain number ("1.00", "0.523", ...).
# Excel's COM layer auto-coerces such strings into numeric cells unless
# the cell is explicitly formatted as Text first; ClearFormats() afterwards
# drops the temporary "@" number format again so no stray style lingers on
# the cell, matching the original (unstyled) inline-string cells.
$textLookingNumbers = [ordered]@{
    'D4' = '1.00'
    'D5' = '305.73'
    'D6' = '94.96'
    'D7' = '0.568'
    'D9' = '0.523'
    'D10' = '34.47'
    'D11' = '0.0806'
    'D12' = '7.18'
    'D16' = '0.817'
    'D17' = '13.48'
    'D20' = '12.19'
    'D21' = '6.17'
    'D22' = '64.65'
    'D23' = '237.28'
    'D24' = '2.90'
    'D26' = '1.94'
    'D27' = '9.83'
    'D28' = '2.13'
    'D29' = '36.76'
    'D30' = '20.04'
    'D31' = '5.86'
    'D32' = '153.65'
    'D33' = '0.0805'
    'D34' = '3.29'
    'D36' = '0.109'
    'D38' = '1.77'
    'D39' = '14.92'
    'D40' = '3.35'
    'D41' = '3.77'
    'D42' = '0.0299'
    'D45' = '84.22'
    'D46' = '0.186'
    'D47' = '99.55'
    'D48' = '4.89'
    'D49' = '14.60'
    'D50' = '8.05'
    'D51' = '68.44'
}
foreach ($ref in $textLookingNumbers.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $textLookingNumbers[$ref]
    $cell.ClearFormats()
}

# Remaining cells (coin names, links, and the "xx.xx%" volume strings, plus
# price strings containing more than one "." like "43.978.50") are never
# misread as numbers, so a plain Value assignment is enough.
$plainText = [ordered]@{
    'D2' = '43.978.50'
    'D3' = '2.234.40'
    'E3' = '  -0.87%  '
    'E4' = '  +0.12%  '
    'E5' = '  -4.25%  '
    'E6' = '  -6.66%  '
    'E7' = '  -1.55%  '
    'E8' = '  +0.19%  '
    'E9' = '  -5.28%  '
    'E10' = '  -7.81%  '
    'E11' = '  -3.05%  '
    'E12' = '  -5.49%  '
    'E13' = '  -2.75%  '
    'D14' = '2.574.67'
    'E14' = '  -0.91%  '
    'D15' = '2.233.38'
    'E15' = '  -1.16%  '
    'E16' = '  -4.64%  '
    'E17' = '  -6.63%  '
    'D18' = '43.864.98'
    'E18' = '  -0.10%  '
    'D19' = '0.0₃0954'
    'E19' = '  -3.13%  '
    'E20' = '  -8.33%  '
    'E21' = '  -5.37%  '
    'E22' = '  -1.65%  '
    'E23' = '  +0.68%  '
    'E24' = '  -7.72%  '
    'E25' = '  +0.36%  '
    'E26' = '  -7.83%  '
    'E27' = '  -3.72%  '
    'E28' = '  -3.09%  '
    'E29' = '  -3.14%  '
    'E30' = '  -0.84%  '
    'E31' = '  -5.29%  '
    'E32' = '  -4.85%  '
    'E33' = '  -5.42%  '
    'E34' = '  +8.61%  '
    'E35' = '  -2.53%  '
    'E36' = '  -5.96%  '
    'E37' = '  -0.56%  '
    'E38' = '  -8.25%  '
    'E39' = '  -11.48%  '
    'E40' = '  -9.95%  '
    'E41' = '  -10.32%  '
    'E42' = '  -5.17%  '
    'E43' = '  +0.17%  '
    'D44' = '1.741.37'
    'E44' = '  -2.85%  '
    'E45' = '  +2.04%  '
    'E46' = '  -6.63%  '
    'E47' = '  -5.06%  '
    'E48' = '  -6.07%  '
    'B49' = 'EnergySwap'
    'C49' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'E49' = '  -0.36%  '
    'B50' = 'FraxShare'
    'C50' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'E50' = '  -3.76%  '
    'B51' = 'ordi'
    'C51' = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
    'E51' = '  -9.33%  '
}
foreach ($ref in $plainText.Keys) {
    $ws.Range($ref).Value = $plainText[$ref]
}
